# Team_PER_1996: fix PER bug - re-point each row's team label to the
# correct team (the original file had the Team column shifted against the
# PER values) and replace the PER column with the corrected per-minute
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Team = "POR"; PER = 12.96923076923077 },
    @{ Row = 3;  Team = "NJN"; PER = 12.05 },
    @{ Row = 4;  Team = "CLE"; PER = 12.84 },
    @{ Row = 5;  Team = "DAL"; PER = 13.775 },
    @{ Row = 6;  Team = "MIA"; PER = 11.24444444444444 },
    @{ Row = 7;  Team = "SEA"; PER = 14.07692307692308 },
    @{ Row = 8;  Team = "ATL"; PER = 13.25 },
    @{ Row = 9;  Team = "MIL"; PER = 11.32 },
    @{ Row = 10; Team = "LAC"; PER = 13.00769230769231 },
    @{ Row = 11; Team = "VAN"; PER = 11.98461538461538 },
    @{ Row = 12; Team = "DET"; PER = 12.41538461538461 },
    @{ Row = 13; Team = "WSB"; PER = 16.92222222222222 },
    @{ Row = 14; Team = "SAS"; PER = 13.88181818181818 },
    @{ Row = 15; Team = "ORL"; PER = 15.15 },
    @{ Row = 16; Team = "UTA"; PER = 13.04285714285714 },
    @{ Row = 17; Team = "HOU"; PER = 12.93333333333333 },
    @{ Row = 18; Team = "DEN"; PER = 12.54615384615385 },
    @{ Row = 19; Team = "LAL"; PER = 14.00714285714286 },
    @{ Row = 20; Team = "GSW"; PER = 13.26363636363636 },
    @{ Row = 21; Team = "IND"; PER = 13.9 },
    @{ Row = 22; Team = "CHI"; PER = 13.27142857142857 },
    @{ Row = 23; Team = "PHI"; PER = 10.84285714285714 },
    @{ Row = 24; Team = "CHH"; PER = 12.31538461538462 },
    @{ Row = 25; Team = "BOS"; PER = 10.75714285714286 },
    @{ Row = 26; Team = "TOR"; PER = 12.23333333333334 },
    @{ Row = 27; Team = "SAC"; PER = 13.3 },
    @{ Row = 28; Team = "PHO"; PER = 15.50714285714286 },
    @{ Row = 29; Team = "NYK"; PER = 15.025 },
    @{ Row = 30; Team = "MIN"; PER = 10.79166666666666 }
)

foreach ($entry in $rows) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Team
    $ws.Cells.Item($entry.Row, 3).Value = $entry.PER
}
